{"js": "// \"json and templates corrections\"\n//\n// The contract template's product table used merge-field placeholders\n// prefixed with \"items.\" (e.g. {{items.product.name}}) but the data\n// actually gets passed under the \"product\" key, so the placeholders are\n// corrected to {{product.name}}, {{product.quantity}}, {{product.price}}\n// and {{product.cost}}. The first two columns of that table are widened\n// slightly (3828 / 2112 twips instead of 4140 / 1800) to better fit the\n// new label text.\n\nconst body = context.document.body;\n\n// --- 1. Resize the first two columns of the pricing table -------------\n// Word JS reports/accepts column widths in points; the OOXML grid is in\n// twips (1 pt = 20 twips), so 3828/20 = 191.4 and 2112/20 = 105.6.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// The pricing table (\"description of ProductS | Quantity | Unit Price ($) |\n// Total cost ($)\") is the third table in the document.\nconst pricingTable = tables.items[2];\npricingTable.getCell(0, 0).columnWidth = 3828 / 20; // 191.4 pt\npricingTable.getCell(0, 1).columnWidth = 2112 / 20; // 105.6 pt\nawait context.sync();\n\n// --- 2. Fix the merge-field placeholders in that table -----------------\nconst replacements = [\n  [\"items.product.name\", \"product.name\"],\n  [\"items.quantity\", \"product.quantity\"],\n  [\"items.product.price\", \"product.price\"],\n  [\"items.cost\", \"product.cost\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# \"json and templates corrections\"\n#\n# The contract template's product table used merge-field placeholders\n# prefixed with \"items.\" (e.g. {{items.product.name}}) but the data\n# actually gets passed under the \"product\" key, so the placeholders are\n# corrected to {{product.name}}, {{product.quantity}}, {{product.price}}\n# and {{product.cost}}. The first two columns of that table are widened\n# slightly (3828 / 2112 twips instead of 4140 / 1800) to better fit the\n# new label text.\n\n$d = $word.ActiveDocument\n\n# --- 1. Resize the first two columns of the pricing table --------------\n# Word reports/accepts column widths in points; the OOXML grid is in\n# twips (1 pt = 20 twips), so 3828/20 = 191.4 and 2112/20 = 105.6.\n# The pricing table (\"description of ProductS | Quantity | Unit Price ($) |\n# Total cost ($)\") is the third table in the document.\n$pricingTable = $d.Tables.Item(3)\n$pricingTable.Columns.Item(1).Width = 3828 / 20\n$pricingTable.Columns.Item(2).Width = 2112 / 20\n\n# --- 2. Fix the merge-field placeholders in that table ------------------\n$replacements = @(\n    @(\"items.product.name\", \"product.name\"),\n    @(\"items.quantity\", \"product.quantity\"),\n    @(\"items.product.price\", \"product.price\"),\n    @(\"items.cost\", \"product.cost\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
